# Update app to include multi-year boxplots, clean foggy days out of data.
#
# The "tower count" log contains one row per calendar date, but several
# dates were never actually surveyed (too foggy / no observations) and
# were left as placeholder rows holding only a date in column A (one of
# them - 2015-07-22 - additionally carried a literal "too foggy " note
# in column G). Those rows are removed so the sheet only contains rows
# that have real observation data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row numbers (1-based, as they exist in the *current* sheet before any
# deletions) that hold nothing but a date - i.e. no squee/herg/gbbg/coei/
# chix/other counts were ever recorded for that day. Row 45 is the one
# "too foggy" day that does carry text in column G; it gets removed too.
# Deleting from the bottom up keeps the remaining row numbers stable as
# we go.
$emptyRows = @(51, 50, 45, 44, 43, 42, 38, 37, 31, 28, 27, 23, 21)

foreach ($r in $emptyRows) {
    $ws.Rows.Item($r).Delete()
}

# Leave the sheet's view matching where the author ended up after the
# cleanup - selecting the last data row/col of the trimmed table.
$ws.Range("D40").Select()
